# Weekly price-sheet update: a new week's record is inserted as row 7
# (pushing the existing rows 7-41 down to rows 8-42), growing the used
# range from A1:R41 to A1:R42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7:41 down one row to make room for the new weekly record.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new week's data.
$ws.Cells.Item(7, 1).Value  = 2
$ws.Cells.Item(7, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(7, 3).Value  = "Coquimbo"
$ws.Cells.Item(7, 4).Value  = 44552
$ws.Cells.Item(7, 5).Value  = 4
$ws.Cells.Item(7, 6).Value  = 100112032
$ws.Cells.Item(7, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(7, 8).Value  = "Sin especificar"
$ws.Cells.Item(7, 9).Value  = "Primera"
$ws.Cells.Item(7, 10).Value = 600
$ws.Cells.Item(7, 11).Value = 6000
$ws.Cells.Item(7, 12).Value = 7000
$ws.Cells.Item(7, 13).Value = 6500
$ws.Cells.Item(7, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 108
$ws.Cells.Item(7, 17).Value = 60
$ws.Cells.Item(7, 18).Value = "Hortaliza"
